$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 708
$ws1.Range("F3").Value = 52
$ws1.Range("F9").Value = 4519
$ws1.Range("F12").Value = 17

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 708
$ws4.Range("F3").Value = 52
$ws4.Range("F9").Value = 4519
$ws4.Range("F12").Value = 17
